# Updated the code so that it can load in a random animal with a random
# stage, with each stage having a 'difficulty' which affects the spawned
# animal. This rewrites the "Difficulty" column (N) on Sheet1 from the old
# 0-100 percentage-style numbers to the new small (signed) difficulty
# integers used by the game logic, and restores Sheet1 (rather than
# Sheet2) as the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New "Difficulty" values for rows 2-19 (column N), keyed by row number.
$difficulty = @{
    2  = -5
    3  = -3
    4  = -2
    5  = 0
    6  = 3
    7  = 2
    8  = 3
    9  = 4
    10 = 2
    11 = 1
    12 = 3
    13 = 4
    14 = 3
    15 = 4
    16 = 3
    17 = 4
    18 = 5
    19 = 5
}

foreach ($row in $difficulty.Keys) {
    $ws1.Range("N$row").Value = $difficulty[$row]
}

# Sheet1 becomes the active sheet/tab again (the source workbook had
# Sheet2 active; the edited workbook has Sheet1 active), with the
# bottom-right frozen pane's selection parked on P15.
[void]$ws1.Activate()
[void]$ws1.Range("P15").Select()
